$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 553.9339164698035
$ws.Range("D2").Value = 136.0796753216788
$ws.Range("F2").Value = 449
$ws.Range("G2").Value = 513
$ws.Range("H2").Value = 625
$ws.Range("C3").Value = 37.54371199562761
$ws.Range("D3").Value = 6.49341519441723
$ws.Range("E3").Value = 13.99
$ws.Range("F3").Value = 32.81
$ws.Range("G3").Value = 37.81
$ws.Range("H3").Value = 41.51
$ws.Range("C4").Value = 1.98216252087201
$ws.Range("D4").Value = 2.542316588933395
$ws.Range("F4").Value = 0.63
$ws.Range("H4").Value = 2.41
$ws.Range("C5").Value = 323.3209024386936
$ws.Range("D5").Value = 10.80527007418962
$ws.Range("F5").Value = 316.54
$ws.Range("G5").Value = 324.78
$ws.Range("H5").Value = 331.49
$ws.Range("C6").Value = 21.20733341915513
$ws.Range("D6").Value = 2.579251734899544
$ws.Range("F6").Value = 19.73
$ws.Range("G6").Value = 21.17
$ws.Range("H6").Value = 22.53
$ws.Range("I6").Value = 40.24
$ws.Range("C7").Value = -76.90280080160562
$ws.Range("D7").Value = 22.89926493012658
$ws.Range("F7").Value = -93
$ws.Range("C8").Value = 7.419475878633727
$ws.Range("D8").Value = 7.098372356541219
$ws.Range("I8").Value = 19
$ws.Range("C9").Value = 9.322680892004572
$ws.Range("D9").Value = 1.685704350766166
$ws.Range("C10").Value = 867.8301139623015
$ws.Range("D10").Value = 0.461504525281789
$ws.Range("C11").Value = 0.5559225975394744
$ws.Range("D11").Value = 0.5890156371236127
$ws.Range("C12").Value = 22.74623852133575
$ws.Range("D12").Value = 12.29406113191247
$ws.Range("C13").Value = 0.6740255086446632
$ws.Range("D13").Value = 0.7506961663181104
$ws.Range("C14").Value = 1.827532081539733
$ws.Range("D14").Value = 1.66442692177414
$ws.Range("C15").Value = 94.1628008016058
$ws.Range("D15").Value = 22.89926493004305
$ws.Range("H15").Value = 110.26
$ws.Range("C16").Value = -86.02480107908279
$ws.Range("D16").Value = 20.47195713110109
$ws.Range("F16").Value = -102.7376019773414
$ws.Range("G16").Value = -84.23249407632485
$ws.Range("C17").Value = -78.60532520044909
$ws.Range("D17").Value = 25.44454038663521
$ws.Range("F17").Value = -93.87736039420676
$ws.Range("G17").Value = -73.57382219273629
